$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 8 and 9 (columns A,B,E,F,G,H,Q,R,Z,AB) ---
$cols8 = @("A","B","E","F","G","H","Q","R","Z","AB")
foreach ($col in $cols8) {
    $r8 = "$col`8"
    $r9 = "$col`9"
    $v8 = $ws.Range($r8).Value2
    $v9 = $ws.Range($r9).Value2
    $ws.Range($r8).Value2 = $v9
    $ws.Range($r9).Value2 = $v8
}

# --- Swap rows 13 and 14 (columns A,Q,R,Z,AB) ---
$cols1314 = @("A","Q","R","Z","AB")
foreach ($col in $cols1314) {
    $r13 = "$col`13"
    $r14 = "$col`14"
    $v13 = $ws.Range($r13).Value2
    $v14 = $ws.Range($r14).Value2
    $ws.Range($r13).Value2 = $v14
    $ws.Range($r14).Value2 = $v13
}

# --- Rotate rows 17,18,19 (columns A,Q,R,Z,AB) ---
# New17 = Old18, New18 = Old19, New19 = Old17
$cols171819 = @("A","Q","R","Z","AB")
foreach ($col in $cols171819) {
    $r17 = "$col`17"
    $r18 = "$col`18"
    $r19 = "$col`19"
    $v17 = $ws.Range($r17).Value2
    $v18 = $ws.Range($r18).Value2
    $v19 = $ws.Range($r19).Value2
    $ws.Range($r17).Value2 = $v18
    $ws.Range($r18).Value2 = $v19
    $ws.Range($r19).Value2 = $v17
}
